# "Generate Report for Handback" - refresh the handback-status report with
# newly observed timestamps for the latest handoff/handback XLIFF generation
# cycle.
#
# Only three cell values actually change (everything else in the diff is
# shared-string table churn caused by the report generator re-serializing
# the sst, not a content change):
#   - zh-cn sheet, row for 39c3eb9d...md:
#       Correspond Handoff Datetime (H2)  2016-08-27 20:46:43 -> 2016-08-27 20:47:36
#       Correspond Handback DateTime (K2) 2016-08-27 20:47:11 -> 2016-08-27 20:47:54
#   - de-de sheet, row for 39c3eb9d...md:
#       Correspond Handback DateTime (K2) 2016-08-27 20:47:18 -> 2016-08-27 20:48:03

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 20:47:36"
$wsZhCn.Range("K2").Value = "2016-08-27 20:47:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-27 20:48:03"
